# Rename "Grades" -> "Levels", add a hidden "Data" sheet with level/level-type
# lookup lists, widen the table + add list data validations driven by that
# sheet, update the J.H.S wording, and tweak a few view/print settings.

$wb = $excel.ActiveWorkbook
$levels = $wb.Worksheets.Item(1)
$levels.Name = "Levels"

# --- Levels sheet: reword the J.H.S rows, resize column A, move selection ---
$levels.Range("A12").Value = "Junior High School 1"
$levels.Range("A13").Value = "Junior High School 2"
$levels.Range("A14").Value = "Junior High School 3"

$levels.Range("A1").ColumnWidth = 25.6640625

# --- Expand the "Grades" table to A1:B100 ---
$tbl = $levels.ListObjects.Item(1)
$tbl.Resize($levels.Range("A1:B100"))

# --- Add the hidden "Data" sheet right after Levels ---
$data = $wb.Worksheets.Add($null, $levels)
$data.Name = "Data"

$levelNames = @(
  "Level",
  "Day Care",
  "Creche",
  "Nursery 1",
  "Nursery 2",
  "Kindergarten 1",
  "Kindergarten 2",
  "Basic 1",
  "Basic 2",
  "Basic 3",
  "Basic 4",
  "Basic 5",
  "Basic 6",
  "Basic 7",
  "Basic 8",
  "Basic 9",
  "Basic 10",
  "Basic 11",
  "Basic 12",
  "Class 1",
  "Class 2",
  "Class 3",
  "Class 4",
  "Class 5",
  "Class 6",
  "Class 7",
  "Class 8",
  "Class 9",
  "Class 10",
  "Class 11",
  "Class 12",
  "Stage 1",
  "Stage 2",
  "Stage 3",
  "Stage 4",
  "Stage 5",
  "Stage 6",
  "Stage 7",
  "Stage 8",
  "Stage 9",
  "Stage 10",
  "Stage 11",
  "Stage 12",
  "Junior High School 1",
  "Junior High School 2",
  "Junior High School 3",
  "Senior High School 1",
  "Senior High School 2",
  "Senior High School 3"
)

$levelTypes = @(
  "Level Type",
  "A","B","C","D","E","F","G","H","I","J","K","L","M",
  "N","O","P","Q","R","S","T","U","V","W","X","Y","Z"
)

for ($i = 0; $i -lt $levelNames.Count; $i++) {
  $data.Cells.Item($i + 1, 1).Value = $levelNames[$i]
}
for ($i = 0; $i -lt $levelTypes.Count; $i++) {
  $data.Cells.Item($i + 1, 2).Value = $levelTypes[$i]
}

$data.Columns.Item(1).ColumnWidth = 17.6640625

# --- Data validation lists on the Levels sheet, sourced from Data ---
$levels.Range("A2:A100").Validation.Add(3, 1, 1, "=Data!`$A`$2:`$A`$49")
$levels.Range("B2:B100").Validation.Add(3, 1, 1, "=Data!`$B`$2:`$B`$27")

# --- Page setup on both sheets ---
$levels.PageSetup.PaperSize = 9
$levels.PageSetup.Orientation = 1
$data.PageSetup.PaperSize = 9
$data.PageSetup.Orientation = 1

# --- Protect the Data sheet and hide it ---
[void]$data.Range("C7").Select()
$data.Protect("")
$data.Visible = $false

# --- Selection / view tweaks ---
[void]$levels.Range("B8").Select()
